$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2260
$ws.Range("I20").Value = 2260
$ws.Range("K20").Value = 2260
$ws.Range("M20").Value = -2030
$ws.Range("H35").Value = 2260
$ws.Range("I35").Value = 2260
$ws.Range("K35").Value = 2260
$ws.Range("M35").Value = -1881
$ws.Range("H62").Value = 2017.3334
$ws.Range("I62").Value = 2161.6667
$ws.Range("J62").Value = 1728.6666
$ws.Range("K62").Value = 2161.6667
$ws.Range("L62").Value = 1728.6666
$ws.Range("M62").Value = -1537.6667
$ws.Range("N62").Value = -2976.6666
$ws.Range("H65").Value = 2017.3334
$ws.Range("I65").Value = 2161.6667
$ws.Range("J65").Value = 1728.6666
$ws.Range("K65").Value = 10808.3335
$ws.Range("L65").Value = 8643.333000000001
$ws.Range("M65").Value = -7688.333500000001
$ws.Range("N65").Value = -14883.333
$ws.Range("H70").Value = 3796.0833
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 3592.1667
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 10776.5001
$ws.Range("M70").Value = -11730
$ws.Range("N70").Value = -11316.5001
$ws.Range("H73").Value = 3796.0833
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 3592.1667
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 10776.5001
$ws.Range("M73").Value = -11064
$ws.Range("N73").Value = -12648.5001
$ws.Range("H113").Value = 3457.9412
$ws.Range("J113").Value = 3564.4443
$ws.Range("L113").Value = 3564.4443
$ws.Range("N113").Value = -10072.4443
$ws.Range("H131").Value = 845
$ws.Range("I131").Value = 638.3333
$ws.Range("J131").Value = 1051.6666
$ws.Range("K131").Value = 1914.9999
$ws.Range("L131").Value = 3154.9998
$ws.Range("M131").Value = 3125.0001
$ws.Range("N131").Value = -13234.9998
$ws.Range("H135").Value = 19238.473
$ws.Range("I135").Value = 23337.592
$ws.Range("K135").Value = 210038.328
$ws.Range("M135").Value = -207503.328
$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200
$ws.Range("H137").Value = 2943071.8
$ws.Range("I137").Value = 4001816
$ws.Range("J137").Value = 2115.111
$ws.Range("K137").Value = 12005448
$ws.Range("L137").Value = 6345.333
$ws.Range("M137").Value = -12002898
$ws.Range("N137").Value = -11445.333
$ws.Range("H141").Value = 726.34784
$ws.Range("I141").Value = 691.1818
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 2073.5454
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = 3106.4546
$ws.Range("N141").Value = -14860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 907.125
$ws.Range("I4").Value = 397
$ws.Range("J4").Value = 1417.25
$ws.Range("K4").Value = 397
$ws.Range("L4").Value = 1417.25
$ws.Range("M4").Value = -281
$ws.Range("N4").Value = -1649.25
$ws.Range("H61").Value = 20449724
$ws.Range("I61").Value = 22750694
$ws.Range("J61").Value = 201173
$ws.Range("K61").Value = 22750694
$ws.Range("L61").Value = 201173
$ws.Range("M61").Value = -22750482
$ws.Range("N61").Value = -201597
$ws.Range("H74").Value = 6633107
$ws.Range("I74").Value = 8656419
$ws.Range("J74").Value = 113546
$ws.Range("K74").Value = 8656419
$ws.Range("L74").Value = 113546
$ws.Range("M74").Value = -8655545
$ws.Range("N74").Value = -115294
$ws.Range("H77").Value = 6633107
$ws.Range("I77").Value = 8656419
$ws.Range("J77").Value = 113546
$ws.Range("K77").Value = 43282095
$ws.Range("L77").Value = 567730
$ws.Range("M77").Value = -43277727
$ws.Range("N77").Value = -576466
$ws.Range("H132").Value = 68985.13
$ws.Range("I132").Value = 48696.81
$ws.Range("J132").Value = 111590.6
$ws.Range("K132").Value = 146090.43
$ws.Range("L132").Value = 334771.8
$ws.Range("M132").Value = -143560.43
$ws.Range("N132").Value = -339831.8
$ws.Range("H136").Value = 20449724
$ws.Range("I136").Value = 22750694
$ws.Range("J136").Value = 201173
$ws.Range("K136").Value = 68252082
$ws.Range("L136").Value = 603519
$ws.Range("M136").Value = -68249532
$ws.Range("N136").Value = -608619

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H134").Value = 6239.25
$ws.Range("I134").Value = 5683.2915
$ws.Range("J134").Value = 9575
$ws.Range("K134").Value = 17049.8745
$ws.Range("L134").Value = 28725
$ws.Range("M134").Value = -14514.8745
$ws.Range("N134").Value = -33795

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3016.7715
$ws.Range("I31").Value = 1421.591
$ws.Range("J31").Value = 5716.3076
$ws.Range("K31").Value = 1421.591
$ws.Range("L31").Value = 5716.3076
$ws.Range("M31").Value = -1126.591
$ws.Range("N31").Value = -6306.3076
$ws.Range("H34").Value = 3016.7715
$ws.Range("I34").Value = 1421.591
$ws.Range("J34").Value = 5716.3076
$ws.Range("K34").Value = 1421.591
$ws.Range("L34").Value = 5716.3076
$ws.Range("M34").Value = -1219.591
$ws.Range("N34").Value = -6120.3076
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2376
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11880
$ws.Range("N65").ClearContents()
$ws.Range("H132").Value = 45363.348
$ws.Range("I132").Value = 30998.146
$ws.Range("J132").Value = 86064.75
$ws.Range("K132").Value = 92994.43799999999
$ws.Range("L132").Value = 258194.25
$ws.Range("M132").Value = -90464.43799999999
$ws.Range("N132").Value = -263254.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 848.75
$ws.Range("I17").Value = 1122.5
$ws.Range("J17").Value = 575
$ws.Range("K17").Value = 3367.5
$ws.Range("L17").Value = 1725
$ws.Range("M17").Value = -3198.5
$ws.Range("N17").Value = -2063
$ws.Range("H25").Value = 800
$ws.Range("I25").Value = 800
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2400
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -2231
$ws.Range("N25").ClearContents()
$ws.Range("H30").Value = 800
$ws.Range("I30").Value = 800
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2400
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2298
$ws.Range("N30").ClearContents()
$ws.Range("H122").Value = 721.7619
$ws.Range("I122").Value = 273
$ws.Range("K122").Value = 2457
$ws.Range("M122").Value = -7
$ws.Range("H131").Value = 13699615
$ws.Range("J131").Value = 1063.3594
$ws.Range("L131").Value = 3190.0782
$ws.Range("N131").Value = -13270.0782

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 997.75
$ws.Range("I102").Value = 860.0833
$ws.Range("J102").Value = 1410.75
$ws.Range("K102").Value = 860.0833
$ws.Range("L102").Value = 1410.75
$ws.Range("M102").Value = 761.9167
$ws.Range("N102").Value = -4654.75
$ws.Range("H132").Value = 50510.438
$ws.Range("I132").Value = 30892.117
$ws.Range("J132").Value = 145799.42
$ws.Range("K132").Value = 92676.351
$ws.Range("L132").Value = 437398.26
$ws.Range("M132").Value = -90146.351
$ws.Range("N132").Value = -442458.26

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 25915.238
$ws.Range("I132").Value = 13146.777
$ws.Range("J132").Value = 57836.39
$ws.Range("K132").Value = 39440.331
$ws.Range("L132").Value = 173509.17
$ws.Range("M132").Value = -36910.331
$ws.Range("N132").Value = -178569.17
$ws.Range("H136").Value = 54987.95
$ws.Range("I136").Value = 34573.305
$ws.Range("K136").Value = 103719.915
$ws.Range("M136").Value = -101169.915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1990.5333
$ws.Range("I122").Value = 1549.2188
$ws.Range("J122").Value = 3076.8462
$ws.Range("K122").Value = 4647.6564
$ws.Range("L122").Value = 9230.5386
$ws.Range("M122").Value = -2197.6564
$ws.Range("N122").Value = -14130.5386
$ws.Range("H129").Value = 32193.334
$ws.Range("J129").Value = 32193.334
$ws.Range("L129").Value = 32193.334
$ws.Range("N129").Value = -42193.334
$ws.Range("H132").Value = 78338.17
$ws.Range("I132").Value = 67151.164
$ws.Range("J132").Value = 111899.2
$ws.Range("K132").Value = 201453.492
$ws.Range("L132").Value = 335697.6
$ws.Range("M132").Value = -198923.492
$ws.Range("N132").Value = -340757.6
$ws.Range("H136").Value = 49205.785
$ws.Range("I136").Value = 39383.42
$ws.Range("J136").Value = 65167.125
$ws.Range("K136").Value = 118150.26
$ws.Range("L136").Value = 195501.375
$ws.Range("M136").Value = -115600.26
$ws.Range("N136").Value = -200601.375
